$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (from 12.5546875 to 13.5546875 in the stored XML).
# This runtime's ColumnWidth (chars) -> stored XML width conversion adds an
# offset of 5/6 and then rounds to the nearest 1/6, so we back that out here
# to land as close as possible to the target stored width of 13.5546875.
$ws.Columns.Item(1).ColumnWidth = (13.5546875 - (5/6))

# Update cell values in A1:A59
$ws.Range("A1").Value = 0.4605107229613179
$ws.Range("A2").Value = 0.749600940611022
$ws.Range("A3").Value = 0.84378233680930959
$ws.Range("A4").Value = 0.78043010293449411
$ws.Range("A5").Value = 0.53786083894927783
$ws.Range("A6").Value = 0.75058911831883479
$ws.Range("A7").Value = 0.82212435736177936
$ws.Range("A8").Value = 0.8948709795594989
$ws.Range("A9").Value = 0.97461329164050958
$ws.Range("A10").Value = 0.83160001445933662
$ws.Range("A11").Value = 0.92482238729982369
$ws.Range("A12").Value = 0.51412592787567379
$ws.Range("A13").Value = 0.80790147634586129
$ws.Range("A14").Value = 0.77310399063661939
$ws.Range("A15").Value = 0.48171721800487566
$ws.Range("A16").Value = 0.9414440318114673
$ws.Range("A17").Value = 0.70555388058520208
$ws.Range("A18").Value = 0.63919335826256973
$ws.Range("A19").Value = 0.72911423685958054
$ws.Range("A20").Value = 0.55590771333753586
$ws.Range("A21").Value = 0.36522072482278184
$ws.Range("A22").Value = 0.13551421085294407
$ws.Range("A23").Value = 0.29986791505640281
$ws.Range("A24").Value = 0.35647071364752653
$ws.Range("A25").Value = 0.47833554133527301
$ws.Range("A26").Value = 0.60960533573143871
$ws.Range("A27").Value = 0.78079410756630996
$ws.Range("A28").Value = 0.56551189956343095
$ws.Range("A29").Value = 0.4557728705613332
$ws.Range("A30").Value = 0.32236526931202641
$ws.Range("A31").Value = 0.28582176328577386
$ws.Range("A32").Value = 0.33595004834293479
$ws.Range("A33").Value = 0.43217313302017379
$ws.Range("A34").Value = 0.31526914221593383
$ws.Range("A35").Value = 0.50016444246607106
$ws.Range("A36").Value = 0.35137368682741782
$ws.Range("A37").Value = 0.48693456171787042
$ws.Range("A38").Value = 0.41616577760197559
$ws.Range("A39").Value = 0.58691629571800108
$ws.Range("A40").Value = 0.50353665321129126
$ws.Range("A41").Value = 0.30697273823764004
$ws.Range("A42").Value = 0.36668374045254981
$ws.Range("A43").Value = 0.21905552587398006
$ws.Range("A44").Value = 0.16797525605635169
$ws.Range("A45").Value = 0.30696802984269594
$ws.Range("A46").Value = 0.15985582878980625
$ws.Range("A47").Value = 0.20225265717246357
$ws.Range("A48").Value = 0.11643630441764816
$ws.Range("A49").Value = 0.023561415710373762
$ws.Range("A50").Value = 0.024174930993493129
$ws.Range("A51").Value = 0.016213125049177539
$ws.Range("A52").Value = 0.017380879015799988
$ws.Range("A53").Value = 0.031178057297968055
$ws.Range("A54").Value = 0.096068992242045223
$ws.Range("A55").Value = 0.13158786456568827
$ws.Range("A56").Value = 0.16672175078823162
$ws.Range("A57").Value = 0.25725101755816759
$ws.Range("A58").Value = 0.13416157424435404
$ws.Range("A59").Value = 0.12111168938950004
